# Daily attendance processing - reorder "Recorded By" (column G) names
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = @($val -split ",\s*")
        if ($parts.Count -gt 1) {
            $reversed = @($parts[-1..-($parts.Count)])
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
